$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) / volume-change (E) data scraped on Sun Aug 13 03:22:20 UTC 2023
$updates = @(
    @{Row=2; D="29.426.71"; E="  +0.15%  "}
    @{Row=3; D="1.850.53"; E="  +0.14%  "}
    @{Row=4; D="1.001"; E="  +0.15%  "}
    @{Row=5; D="240.11"; E="  -0.03%  "}
    @{Row=6; D="0.6290"; E="  -0.16%  "}
    @{Row=7; D=""; E="  +0.12%  "}
    @{Row=8; D="0.07646"; E="  +0.32%  "}
    @{Row=9; D="0.2912"; E="  -0.60%  "}
    @{Row=10; D="24.83"; E="  +1.46%  "}
    @{Row=11; D="2.102.72"; E="  +13.67%  "}
    @{Row=12; D="0.07743"; E="  +0.06%  "}
    @{Row=13; D="5.033"; E="  +0.61%  "}
    @{Row=14; D="0.6807"; E="  +0.27%  "}
    @{Row=15; D=""; E="  -4.79%  "}
    @{Row=16; D="83.48"; E="  -0.21%  "}
    @{Row=17; D="6.186"; E="  +0.09%  "}
    @{Row=18; D="29.521.08"; E="  +0.32%  "}
    @{Row=19; D="228.65"; E="  +0.01%  "}
    @{Row=20; D="12.33"; E="  -0.87%  "}
    @{Row=21; D="1.001"; E="  +0.09%  "}
    @{Row=22; D="7.543"; E="  +0.80%  "}
    @{Row=23; D=""; E="  +0.14%  "}
    @{Row=24; D="157.51"; E="  +0.14%  "}
    @{Row=25; D="0.1385"; E="  -0.83%  "}
    @{Row=26; D="8.426"; E="  +0.93%  "}
    @{Row=27; D="17.72"; E="  +0.62%  "}
    @{Row=28; D="1.399"; E="  +7.70%  "}
    @{Row=29; D="1.463"; E="  +0.07%  "}
    @{Row=30; D="0.05606"; E="  +0.29%  "}
    @{Row=31; D="4.130"; E="  +0.51%  "}
    @{Row=32; D="4.051"; E="  +0.51%  "}
    @{Row=33; D="1.846"; E="  +0.01%  "}
    @{Row=34; D="1.164"; E="  +0.61%  "}
    @{Row=35; D="0.6955"; E="  -2.08%  "}
    @{Row=36; D="2.588"; E="  +0.23%  "}
    @{Row=37; D="0.01802"; E="  -0.04%  "}
    @{Row=38; D="1.230.55"; E="  -0.75%  "}
    @{Row=39; D=""; E="  -1.10%  "}
    @{Row=40; D="6.448"; E="  +0.50%  "}
    @{Row=41; D="0.9098"; E="  +0.52%  "}
    @{Row=42; D="1.001"; E="  +0.10%  "}
    @{Row=43; D="102.58"; E="  +0.69%  "}
    @{Row=44; D="66.04"; E="  +0.26%  "}
    @{Row=45; D="7.195"; E="  +0.52%  "}
    @{Row=46; D="0.00000000118"; E="  -2.86%  "}
    @{Row=47; D="0.4029"; E="  +0.28%  "}
    @{Row=48; D="9.021"; E="  -0.19%  "}
    @{Row=49; D=""; E="  +2.96%  "}
    @{Row=50; D="1.684"; E="  +0.17%  "}
    @{Row=51; D="0.05708"; E="  +0.01%  "}
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.D -ne "") {
        $dCell = $ws.Cells.Item($row, 4)
        $dVal = $u.D

        # Values that look numeric (a single decimal point, parseable as a
        # plain number) get silently coerced/rounded by Excel's normal
        # type-inference when assigned through .Value (e.g. "1.001" ->
        # 1.0009999999999999, or trailing zeros get stripped). Force the
        # cell to Text first so the literal digit string round-trips
        # exactly, matching the source data export. Values that already
        # contain two decimal points (e.g. "29.426.71") are never valid
        # numbers, so Excel keeps them as text on its own.
        $isPlainNumber = $dVal -match "^[+-]?[0-9]*\.?[0-9]+$"

        if ($isPlainNumber) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $dVal
    }

    if ($u.E -ne "") {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
